{"js": "// Replace the cover-letter body paragraphs per the diff.\nconst body = context.document.body;\n\n// 1) Drop the entire \"keywords\" paragraph (it has no counterpart in the new text),\n//    together with the pair of <w:br/> line breaks that trail it.\nconst kwHits = body.search(\"The keywords in your job description, such as \\\"Software Development,\\\" \\\"Python,\\\" \\\"C++,\\\" and \\\"Data-Driven Approach,\\\" align strongly with my skillset. My proficiency in these technologies, coupled with my passion for building innovative solutions, make me an ideal candidate for this role.\", {matchCase: true});\nconst eagerHits = body.search(\"I am eager to contribute my skills to your cutting-edge software development projects. I am confident that my dedication to delivering high-quality work, collaborative spirit, and drive to learn and grow will enable me to thrive in your dynamic work environment.\", {matchCase: true});\nkwHits.load('items');\neagerHits.load('items');\nawait context.sync();\nif (kwHits.items.length && eagerHits.items.length) {\n    const kwRange = kwHits.items[0];\n    const eagerStart = eagerHits.items[0].getRange('Start');\n    const toRemove = kwRange.expandTo(eagerStart);\n    toRemove.delete();\n    await context.sync();\n}\n\n// 2) Swap the remaining paragraph bodies for their rewritten versions, one at a time,\n//    so each run keeps its original formatting (rPr) and surrounding <w:br/> breaks.\nconst replacements = [\n    [\"I am writing to express my keen interest in the Software Developer Intern position at your esteemed company. With my strong academic background in computer science, relevant work experience, and eagerness to contribute to innovative software development projects, I am confident that I can make a significant contribution to your team.\", \"I am writing to express my interest in the Software Developer Intern position at your company. I am currently pursuing a Master's degree in Computer Science from Arizona State University with a focus on Software Development and am expected to graduate in December 2025.\"],\n    [\"During my academic journey, I have acquired a solid foundation in software development principles, algorithms, and data structures. My coursework in Data Visualisation, Cloud Computing, and Frontier Topics in GenAI has equipped me with the necessary skills to tackle complex technical challenges in the industry. Moreover, my research experience in statistical analytics has honed my ability to extract meaningful insights from large datasets.\", \"Throughout my academic and professional experience, I have developed a strong foundation in software development principles, including object-oriented design, data structures, and algorithms. I am proficient in programming languages such as Python, Java, and C++, and I have experience working with various software development tools and technologies.\"],\n    [\"My professional experience at FIS and OG Advertising Private Limited further solidified my technical capabilities. As a Software Engineer at FIS, I developed and deployed a RPA bot using Selenium and Java, resulting in significant cost savings for the company. Through my work at OG Advertising, I gained expertise in web development and implemented a custom CRM system that streamlined admission processes for educational institutes.\", \"In my previous role as a Software Engineer at FIS, I was responsible for designing and developing an RPA bot using Selenium and Java. This bot automated complex data operations, resulting in annual savings of 2100 hours. I also enhanced the performance of two existing RPA bots, reducing runtime by 20% and improving overall efficiency.\"],\n    [\"I am eager to contribute my skills to your cutting-edge software development projects. I am confident that my dedication to delivering high-quality work, collaborative spirit, and drive to learn and grow will enable me to thrive in your dynamic work environment.\", \"I am confident that my skills and experience align well with the requirements of the Software Developer Intern position. I am eager to contribute to your team and learn from experienced professionals in the field. I am a highly motivated and results-oriented individual with a strong work ethic and excellent communication skills.\"],\n    [\"Thank you for your time and consideration. I look forward to the opportunity to further discuss my qualifications and how I can contribute to the success of your team.\", \"Thank you for your time and consideration. I look forward to the opportunity to further discuss my qualifications and how I can contribute to the success of your company.\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n    const hits = body.search(oldText, {matchCase: true});\n    hits.load('items');\n    await context.sync();\n    if (!hits.items.length) {\n        throw new Error('Expected text not found: ' + oldText.substring(0, 60));\n    }\n    hits.items[0].insertText(newText, 'Replace');\n    await context.sync();\n}\n", "ps1": "# Rewrite the cover-letter body paragraphs per the target diff.\n# The whole letter body lives as <w:br/>-separated <w:t> runs inside ONE paragraph/run,\n# so paragraph-level Find/Replace would touch the whole run; instead we locate each exact\n# old sentence by character offset (via .IndexOf on the live text) and rewrite just that\n# span with Range.Text, which leaves the surrounding <w:br/> breaks and run formatting\n# untouched and - unlike Find.Execute's Replace - does not \"smart quote\" straight quotes\n# / apostrophes in the incoming text.\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($oldText, $newText) {\n    $body = $d.Content\n    $bodyText = $body.Text\n    $idx = $bodyText.IndexOf($oldText)\n    if ($idx -lt 0) { throw \"Expected text not found: $($oldText.Substring(0, [Math]::Min(60, $oldText.Length)))\" }\n    $start = $body.Start + $idx\n    $end = $start + $oldText.Length\n    $d.Range($start, $end).Text = $newText\n}\n\n# 1) Drop the entire \"keywords\" paragraph (it has no counterpart in the new text),\n#    together with the pair of line-break characters that trail it.\n$kwText = 'The keywords in your job description, such as \"Software Development,\" \"Python,\" \"C++,\" and \"Data-Driven Approach,\" align strongly with my skillset. My proficiency in these technologies, coupled with my passion for building innovative solutions, make me an ideal candidate for this role.'\n$body = $d.Content\n$bodyText = $body.Text\n$kwIdx = $bodyText.IndexOf($kwText)\nif ($kwIdx -lt 0) { throw \"keywords paragraph not found\" }\n$delStart = $body.Start + $kwIdx\n$delEnd = $delStart + $kwText.Length + 2  # swallow the following <w:br/><w:br/> pair\n$d.Range($delStart, $delEnd).Delete()\n\n# 2) Swap the remaining paragraph bodies for their rewritten versions, one at a time.\nReplace-ExactText 'I am writing to express my keen interest in the Software Developer Intern position at your esteemed company. With my strong academic background in computer science, relevant work experience, and eagerness to contribute to innovative software development projects, I am confident that I can make a significant contribution to your team.' 'I am writing to express my interest in the Software Developer Intern position at your company. I am currently pursuing a Master''s degree in Computer Science from Arizona State University with a focus on Software Development and am expected to graduate in December 2025.'\nReplace-ExactText 'During my academic journey, I have acquired a solid foundation in software development principles, algorithms, and data structures. My coursework in Data Visualisation, Cloud Computing, and Frontier Topics in GenAI has equipped me with the necessary skills to tackle complex technical challenges in the industry. Moreover, my research experience in statistical analytics has honed my ability to extract meaningful insights from large datasets.' 'Throughout my academic and professional experience, I have developed a strong foundation in software development principles, including object-oriented design, data structures, and algorithms. I am proficient in programming languages such as Python, Java, and C++, and I have experience working with various software development tools and technologies.'\nReplace-ExactText 'My professional experience at FIS and OG Advertising Private Limited further solidified my technical capabilities. As a Software Engineer at FIS, I developed and deployed a RPA bot using Selenium and Java, resulting in significant cost savings for the company. Through my work at OG Advertising, I gained expertise in web development and implemented a custom CRM system that streamlined admission processes for educational institutes.' 'In my previous role as a Software Engineer at FIS, I was responsible for designing and developing an RPA bot using Selenium and Java. This bot automated complex data operations, resulting in annual savings of 2100 hours. I also enhanced the performance of two existing RPA bots, reducing runtime by 20% and improving overall efficiency.'\nReplace-ExactText 'I am eager to contribute my skills to your cutting-edge software development projects. I am confident that my dedication to delivering high-quality work, collaborative spirit, and drive to learn and grow will enable me to thrive in your dynamic work environment.' 'I am confident that my skills and experience align well with the requirements of the Software Developer Intern position. I am eager to contribute to your team and learn from experienced professionals in the field. I am a highly motivated and results-oriented individual with a strong work ethic and excellent communication skills.'\nReplace-ExactText 'Thank you for your time and consideration. I look forward to the opportunity to further discuss my qualifications and how I can contribute to the success of your team.' 'Thank you for your time and consideration. I look forward to the opportunity to further discuss my qualifications and how I can contribute to the success of your company.'\n"}
